$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from 12 to 17
# (ColumnWidth round-trips through Excel's max-digit-width units, adding a
#  constant ~5/6 offset on save; back the input off so the stored <col>
#  width lands on exactly 17, matching the target diff.)
$ws.Columns.Item(1).ColumnWidth = 17 - 5/6

# Rename processo_2 -> copy_processo_2 in A2
$ws.Range("A2").Value = "copy_processo_2"

# Convert E3 from text "44432" to a real number 44432
$ws.Range("E3").Value = 44432
